$d = $word.ActiveDocument

$d.Content.Find.Execute("812×5=", $true, $false, $false, $false, $false, $true, 1, $false, "978×9=", 2)
$d.Content.Find.Execute("454×6=", $true, $false, $false, $false, $false, $true, 1, $false, "725×6=", 2)
$d.Content.Find.Execute("344×2=", $true, $false, $false, $false, $false, $true, 1, $false, "220×2=", 2)
$d.Content.Find.Execute("126×5=", $true, $false, $false, $false, $false, $true, 1, $false, "558×3=", 2)
$d.Content.Find.Execute("285×3=", $true, $false, $false, $false, $false, $true, 1, $false, "646×4=", 2)
$d.Content.Find.Execute("414×9=", $true, $false, $false, $false, $false, $true, 1, $false, "352×6=", 2)
$d.Content.Find.Execute("377×7=", $true, $false, $false, $false, $false, $true, 1, $false, "818×6=", 2)
$d.Content.Find.Execute("845×9=", $true, $false, $false, $false, $false, $true, 1, $false, "121×8=", 2)
$d.Content.Find.Execute("823×5=", $true, $false, $false, $false, $false, $true, 1, $false, "722×7=", 2)
$d.Content.Find.Execute("638×8=", $true, $false, $false, $false, $false, $true, 1, $false, "936×3=", 2)
$d.Content.Find.Execute("544×3=", $true, $false, $false, $false, $false, $true, 1, $false, "671×4=", 2)
$d.Content.Find.Execute("334×7=", $true, $false, $false, $false, $false, $true, 1, $false, "889×5=", 2)
$d.Content.Find.Execute("788×4=", $true, $false, $false, $false, $false, $true, 1, $false, "489×7=", 2)
$d.Content.Find.Execute("585×8=", $true, $false, $false, $false, $false, $true, 1, $false, "834×6=", 2)
$d.Content.Find.Execute("908×6=", $true, $false, $false, $false, $false, $true, 1, $false, "408×9=", 2)
$d.Content.Find.Execute("809×2=", $true, $false, $false, $false, $false, $true, 1, $false, "907×6=", 2)
$d.Content.Find.Execute("334×3=", $true, $false, $false, $false, $false, $true, 1, $false, "912×6=", 2)
$d.Content.Find.Execute("348×9=", $true, $false, $false, $false, $false, $true, 1, $false, "860×7=", 2)
$d.Content.Find.Execute("683×2=", $true, $false, $false, $false, $false, $true, 1, $false, "140×6=", 2)
$d.Content.Find.Execute("208×3=", $true, $false, $false, $false, $false, $true, 1, $false, "873×3=", 2)
$d.Content.Find.Execute("898×6=", $true, $false, $false, $false, $false, $true, 1, $false, "991×4=", 2)
$d.Content.Find.Execute("353×4=", $true, $false, $false, $false, $false, $true, 1, $false, "482×2=", 2)
$d.Content.Find.Execute("945×8=", $true, $false, $false, $false, $false, $true, 1, $false, "591×9=", 2)
$d.Content.Find.Execute("126×6=", $true, $false, $false, $false, $false, $true, 1, $false, "718×3=", 2)
$d.Content.Find.Execute("197×9=", $true, $false, $false, $false, $false, $true, 1, $false, "226×7=", 2)
